# Update "想去人数" (F column) counts across sheets per the source diff.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 200
$ws1.Range("F5").Value  = 979
$ws1.Range("F6").Value  = 5414
$ws1.Range("F8").Value  = 661
$ws1.Range("F9").Value  = 939
$ws1.Range("F12").Value = 34
$ws1.Range("F14").Value = 24
$ws1.Range("F17").Value = 1807
$ws1.Range("F18").Value = 1461
$ws1.Range("F19").Value = 872
$ws1.Range("F23").Value = 533
$ws1.Range("F28").Value = 2772
$ws1.Range("F32").Value = 115
$ws1.Range("F33").Value = 32
$ws1.Range("F34").Value = 352
$ws1.Range("F40").Value = 679
$ws1.Range("F41").Value = 85
$ws1.Range("F42").Value = 51
$ws1.Range("F43").Value = 53

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 176
$ws2.Range("F6").Value = 126

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 200
$ws4.Range("F5").Value  = 979
$ws4.Range("F7").Value  = 5414
$ws4.Range("F9").Value  = 661
$ws4.Range("F11").Value = 176
$ws4.Range("F12").Value = 939
$ws4.Range("F15").Value = 126
$ws4.Range("F17").Value = 34
$ws4.Range("F19").Value = 24
$ws4.Range("F23").Value = 1807
$ws4.Range("F24").Value = 1461
$ws4.Range("F25").Value = 872
$ws4.Range("F29").Value = 533
$ws4.Range("F32").Value = 2772
$ws4.Range("F36").Value = 115
$ws4.Range("F37").Value = 32
$ws4.Range("F38").Value = 352
$ws4.Range("F43").Value = 679
$ws4.Range("F44").Value = 85
$ws4.Range("F45").Value = 53

$wb.Save()
